# Add "2022-Q1" sheet (fund holdings detail) before the "总计" summary sheet,
# and update "总计" with the new quarter row.

$wb = $excel.ActiveWorkbook

# --- Step 1: create the new "2022-Q1" sheet by duplicating "2021-Q4" ---
# (this keeps identical column widths / header text / cell styles, so we
# only need to overwrite the data rows afterwards)
$templateSheet = $wb.Worksheets.Item("2021-Q4")
$templateSheet.Copy($wb.Worksheets.Item("总计"))
$newSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newSheet.Name = "2022-Q1"

# NOTE: re-fetch "总计" by name *after* the Copy() above — the variable
# captured beforehand would now resolve to the newly inserted sheet
# since it shifted into that sheet's old position.
$totalSheet = $wb.Worksheets.Item("总计")

# the template ("2021-Q4") has 75 data rows (rows 2-76); the new quarter
# only has 40 funds (rows 2-41), so drop the extra rows entirely
$newSheet.Range("A42:H76").Clear()

$fundData = @(
    ,@(0, "003961", "易方达瑞程灵活配置混合A", "34.99", "91.05", "7.68", "2.6872", 2)
    ,@(1, "002168", "嘉实智能汽车股票", "52.15", "92.00", "4.14", "2.1590", 9)
    ,@(2, "001182", "易方达安心回馈混合", "90.36", "34.92", "1.68", "1.5180", 4)
    ,@(3, "470098", "汇添富逆向投资混合", "30.96", "93.72", "4.84", "1.4985", 7)
    ,@(4, "005802", "汇添富智能制造股票", "31.46", "91.61", "4.58", "1.4409", 8)
    ,@(5, "070010", "嘉实主题精选混合", "26.54", "83.50", "4.20", "1.1147", 5)
    ,@(6, "003962", "易方达瑞程灵活配置混合C", "9.83", "91.05", "7.68", "0.7549", 2)
    ,@(7, "001039", "嘉实先进制造股票", "14.88", "93.31", "4.94", "0.7351", 7)
    ,@(8, "000986", "太平灵活配置混合型发起式", "18.13", "82.45", "3.48", "0.6309", 9)
    ,@(9, "008949", "平安匠心优选混合A", "12.05", "89.95", "3.82", "0.4603", 6)
    ,@(10, "004477", "嘉实沪港深回报混合", "8.77", "84.42", "4.28", "0.3754", 6)
    ,@(11, "070018", "嘉实回报混合", "7.73", "74.29", "3.61", "0.2791", 6)
    ,@(12, "003839", "易方达瑞通灵活配置混合A", "9.33", "32.04", "1.86", "0.1735", 3)
    ,@(13, "009064", "鹏扬景沃六个月持有期混合A", "23.92", "21.80", "0.67", "0.1603", 5)
    ,@(14, "001443", "易方达瑞选灵活配置混合I", "11.70", "23.79", "0.99", "0.1158", 2)
    ,@(15, "001444", "易方达瑞选灵活配置混合E", "11.70", "23.79", "0.99", "0.1158", 2)
    ,@(16, "001541", "汇添富民营新动力股票", "2.30", "89.06", "4.60", "0.1058", 4)
    ,@(17, "519668", "银河竞争优势成长混合", "2.12", "86.14", "4.94", "0.1047", 6)
    ,@(18, "009215", "易方达瑞川灵活配置混合A", "10.76", "21.17", "0.97", "0.1044", 4)
    ,@(19, "001314", "易方达新益灵活配置混合I", "9.93", "24.13", "1.05", "0.1043", 2)
    ,@(20, "001315", "易方达新益灵活配置混合E", "9.93", "24.13", "1.05", "0.1043", 2)
    ,@(21, "003882", "易方达瑞弘灵活配置混合A", "6.90", "28.58", "1.50", "0.1035", 5)
    ,@(22, "001747", "易方达瑞祺灵活配置混合I", "9.77", "23.89", "1.01", "0.0987", 3)
    ,@(23, "001748", "易方达瑞祺灵活配置混合E", "9.77", "23.89", "1.01", "0.0987", 3)
    ,@(24, "001688", "嘉实新起点灵活配置混合A", "6.65", "23.30", "1.06", "0.0705", 9)
    ,@(25, "002222", "嘉实新趋势灵活配置混合", "5.84", "23.37", "1.01", "0.0590", 9)
    ,@(26, "001755", "嘉实新思路灵活配置混合", "5.31", "25.63", "1.01", "0.0536", 9)
    ,@(27, "009065", "鹏扬景沃六个月持有期混合C", "7.70", "21.80", "0.67", "0.0516", 5)
    ,@(28, "008950", "平安匠心优选混合C", "1.15", "89.95", "3.82", "0.0439", 6)
    ,@(29, "009537", "太平行业优选股票A", "0.88", "90.50", "3.80", "0.0334", 10)
    ,@(30, "003883", "易方达瑞弘灵活配置混合C", "2.02", "28.58", "1.50", "0.0303", 5)
    ,@(31, "003840", "易方达瑞通灵活配置混合C", "1.38", "32.04", "1.86", "0.0257", 3)
    ,@(32, "001972", "前海开源沪港深智慧生活优选灵活配置混合", "0.48", "77.26", "4.37", "0.0210", 8)
    ,@(33, "080002", "长盛创新先锋混合", "0.83", "74.24", "2.52", "0.0209", 10)
    ,@(34, "002149", "嘉实新优选灵活配置混合", "0.22", "93.76", "7.97", "0.0175", 3)
    ,@(35, "001864", "中海魅力长三角灵活配置混合", "0.22", "74.64", "5.00", "0.0110", 1)
    ,@(36, "009538", "太平行业优选股票C", "0.20", "90.50", "3.80", "0.0076", 10)
    ,@(37, "006603", "嘉实互融精选股票", "0.14", "85.91", "3.63", "0.0051", 10)
    ,@(38, "002178", "嘉实新起点灵活配置混合C", "0.26", "23.30", "1.06", "0.0028", 9)
    ,@(39, "009216", "易方达瑞川灵活配置混合C", "0.21", "21.17", "0.97", "0.0020", 4)

)

foreach ($row in $fundData) {
    $r = [int]$row[0] + 2
    $newSheet.Cells.Item($r, 1).Value = [int]$row[0]

    $newSheet.Cells.Item($r, 2).Value = "'" + $row[1]
    $newSheet.Cells.Item($r, 2).Style = "Normal"

    $newSheet.Cells.Item($r, 3).Value = $row[2]

    $newSheet.Cells.Item($r, 4).Value = "'" + $row[3]
    $newSheet.Cells.Item($r, 4).Style = "Normal"

    $newSheet.Cells.Item($r, 5).Value = "'" + $row[4]
    $newSheet.Cells.Item($r, 5).Style = "Normal"

    $newSheet.Cells.Item($r, 6).Value = "'" + $row[5]
    $newSheet.Cells.Item($r, 6).Style = "Normal"

    $newSheet.Cells.Item($r, 7).Value = "'" + $row[6]
    $newSheet.Cells.Item($r, 7).Style = "Normal"

    $newSheet.Cells.Item($r, 8).Value = [int]$row[7]
}

# --- Step 2: update the "总计" summary sheet with the new quarter ---
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Cells.Item(2, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value = 40
$totalSheet.Cells.Item(2, 4).Value = 15.5

# renumber the index column (A) for every data row, now that a row was
# inserted at the top (0, 1, 2, 3, 4, 5)
for ($i = 0; $i -le 5; $i++) {
    $r = $i + 2
    $totalSheet.Cells.Item($r, 1).Value = $i
}
